$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be read as text so values such as "318.61"
# are not silently coerced to numbers by Excel's type inference.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '41.604.77'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '2.471.64'
$ws.Range('E3').Value = '  -0.34%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '318.61'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').Value = '92.06'
$ws.Range('E6').Value = '  -0.67%  '
$ws.Range('E7').Value = '  +1.13%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('E9').Value = '  +1.32%  '
$ws.Range('E10').Value = '  +9.41%  '
$ws.Range('D11').Value = '32.86'
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').Value = '2.851.80'
$ws.Range('E13').Value = '  -0.30%  '
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').Value = '15.55'
$ws.Range('E15').Value = '  -4.54%  '
$ws.Range('D16').Value = '2.464.38'
$ws.Range('E16').Value = '  -1.30%  '
$ws.Range('D17').Value = '0.787'
$ws.Range('E17').Value = '  +2.06%  '
$ws.Range('D18').Value = '41.553.25'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').Value = '0.0₃0948'
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('D20').Value = '6.44'
$ws.Range('E20').Value = '  -2.31%  '
$ws.Range('D21').Value = '71.17'
$ws.Range('E21').Value = '  -1.64%  '
$ws.Range('D22').Value = '11.31'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('D23').Value = '238.84'
$ws.Range('E23').Value = '  +1.10%  '
$ws.Range('D24').Value = '2.75'
$ws.Range('E24').Value = '  +0.89%  '
$ws.Range('D25').Value = '1.93'
$ws.Range('E25').Value = '  +1.34%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = '24.62'
$ws.Range('E27').Value = '  -0.76%  '
$ws.Range('E28').Value = '  +3.07%  '
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('D30').Value = '36.20'
$ws.Range('E30').Value = '  +0.96%  '
$ws.Range('D31').Value = '160.91'
$ws.Range('E31').Value = '  +1.67%  '
$ws.Range('D32').Value = '5.50'
$ws.Range('E32').Value = '  +0.93%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('D34').Value = '2.59'
$ws.Range('E34').Value = '  +0.42%  '
$ws.Range('D35').Value = '0.0767'
$ws.Range('E35').Value = '  +1.39%  '
$ws.Range('D36').Value = '17.19'
$ws.Range('E36').Value = '  -2.02%  '
$ws.Range('D37').Value = '2.90'
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('E38').Value = '  +1.29%  '
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('D40').Value = '0.103'
$ws.Range('E40').Value = '  -3.01%  '
$ws.Range('D41').Value = '3.96'
$ws.Range('E41').Value = '  -3.14%  '
$ws.Range('E42').Value = '  +2.98%  '
$ws.Range('D43').Value = '1.990.29'
$ws.Range('E43').Value = '  +1.47%  '
$ws.Range('D44').Value = '0.0285'
$ws.Range('D45').Value = '18.95'
$ws.Range('E45').Value = '  -2.94%  '
$ws.Range('E46').Value = '  +0.87%  '
$ws.Range('E47').Value = '  +2.28%  '
$ws.Range('D48').Value = '2.711.34'
$ws.Range('D49').Value = '97.57'
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('D50').Value = '73.70'
$ws.Range('E50').Value = '  +1.87%  '
$ws.Range('D51').Value = '66.84'
$ws.Range('E51').Value = '  -1.90%  '

# Restore original (default) cell formatting now that the text is committed.
$ws.Range('D2:D51').ClearFormats()
